$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$boldLabel = "Meta description"
$restOfLine = ": Discover the respins mechanic and intuitive gameplay of Break da Bank Again Respins by Microgaming. Play free and read our review to learn more."

$metaStart = $metaPara.Range.Start
$metaPara.Range.Text = $boldLabel + $restOfLine

$labelRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$labelRange.Bold = 1

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph near the end of the document
#    ("Play Break da Bank Again Respins Free | Slot Review"), and replace the
#    text of the paragraph that follows it (the italic meta-description-like
#    paragraph) with the new image-prompt text, keeping its italic run
#    formatting.
# ---------------------------------------------------------------------------
$target = "Play Break da Bank Again Respins Free | Slot Review"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $target) {
        $para.Range.Delete()
        break
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastRange = $d.Range($lastStart, $lastEnd - 1)

$newPromptText = 'Prompt: Create a cartoon-style feature image for "Break da Bank Again Respins" featuring a happy Maya warrior with glasses. Description: The image should show a Maya warrior wearing glasses and a big smile, holding a bag of gold coins in one hand and a slot machine lever in the other. The background should be filled with colorful banknotes and stacks of gold bars. The Maya warrior should be dressed in traditional clothing, with a feather headdress and colorful patterns on his garment. The overall style of the image should be cartoonish and fun, with bright colors and simple shapes. The image should be eye-catching and convey the excitement and joy of a big win on "Break da Bank Again Respins".'

$lastRange.Text = $newPromptText
